# number-showcase.xlsx edit script
# Implements:
#  - inserts a new "sms" command-family column into the hidden '#system' sheet
#    (single command: sendText(phones,text))
#  - inserts a new "ws.async" command-family column into the hidden '#system'
#    sheet (async flavours of the existing web-service commands: download,
#    get, head, patch, post, put)
#  - updates the sorted "target" list (column A) to include the two new
#    command families in their correct alphabetical slots
#  - updates every defined name whose column shifted because of the two
#    new columns, and adds the two new defined names (sms, ws.async)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1. Make room for the two new command-family columns.
#    Before:  ... P(redis) Q(sound) R(ssh) S(step) T(web) U(webalert)
#              V(webcookie) W(ws) X(xml)
#    Inserting at Q pushes sound..xml one column to the right (R..Y).
#    Inserting again at Y (the new position of the old X/xml column)
#    pushes that single column one further to the right, to Z, and frees
#    up Y for the new ws.async family.
# ---------------------------------------------------------------------
$ws.Columns("Q").Insert()
$ws.Columns("Y").Insert()

# ---------------------------------------------------------------------
# 2. Populate the new "sms" column (Q) - one command.
# ---------------------------------------------------------------------
$ws.Range("Q1").Value = "sms"
$ws.Range("Q2").Value = "sendText(phones,text)"

# ---------------------------------------------------------------------
# 3. Populate the new "ws.async" column (Y) - async flavours of the
#    existing web-service commands.
# ---------------------------------------------------------------------
$ws.Range("Y1").Value = "ws.async"
$ws.Range("Y2").Value = "download(url,queryString,saveTo)"
$ws.Range("Y3").Value = "get(url,queryString,output)"
$ws.Range("Y4").Value = "head(url,output)"
$ws.Range("Y5").Value = "patch(url,body,output)"
$ws.Range("Y6").Value = "post(url,body,output)"
$ws.Range("Y7").Value = "put(url,body,output)"

# ---------------------------------------------------------------------
# 4. Re-write the sorted "target" list (column A) so it includes the two
#    new command families ("sms" before "sound", "ws.async" between "ws"
#    and "xml").
# ---------------------------------------------------------------------
$targets = @(
    "target",
    "aws.s3",
    "base",
    "csv",
    "desktop",
    "excel",
    "external",
    "image",
    "io",
    "jms",
    "json",
    "mail",
    "number",
    "pdf",
    "rdbms",
    "redis",
    "sms",
    "sound",
    "ssh",
    "step",
    "web",
    "webalert",
    "webcookie",
    "ws",
    "ws.async",
    "xml"
)
for ($i = 0; $i -lt $targets.Length; $i++) {
    $ws.Range("A" + ($i + 1)).Value = $targets[$i]
}

# ---------------------------------------------------------------------
# 5. Fix up the defined names that refer to columns which shifted right
#    because of the two inserted columns, and add the two brand-new
#    names.
# ---------------------------------------------------------------------
$wb.Names.Item("target").RefersTo      = "='#system'!`$A`$2:`$A`$26"
$wb.Names.Item("sound").RefersTo       = "='#system'!`$R`$2:`$R`$5"
$wb.Names.Item("ssh").RefersTo         = "='#system'!`$S`$2:`$S`$9"
$wb.Names.Item("step").RefersTo        = "='#system'!`$T`$2:`$T`$4"
$wb.Names.Item("web").RefersTo         = "='#system'!`$U`$2:`$U`$108"
$wb.Names.Item("webalert").RefersTo    = "='#system'!`$V`$2:`$V`$6"
$wb.Names.Item("webcookie").RefersTo   = "='#system'!`$W`$2:`$W`$8"
$wb.Names.Item("ws").RefersTo          = "='#system'!`$X`$2:`$X`$16"
$wb.Names.Item("xml").RefersTo         = "='#system'!`$Z`$2:`$Z`$11"

$wb.Names.Add("sms", "='#system'!`$Q`$2:`$Q`$2")
$wb.Names.Add("ws.async", "='#system'!`$Y`$2:`$Y`$7")

# ---------------------------------------------------------------------
# 6. Nudge the sheet's recorded used-range out to AA108 (matches the
#    canonical file's <dimension ref="A1:AA108"/>) without leaving any
#    visible/formatted content behind.
# ---------------------------------------------------------------------
$ws.Range("AA108").Font.Bold = $true
$ws.Range("AA108").Font.Bold = $false
